# Weekly update: insert the new week's price records (2023-08-28, serial
# 45166) for "Agrícola del Norte S.A. de Arica" / Frutilla above the
# existing historical rows 108-110 (dated 2021-02-15, serial 44242),
# pushing them down to rows 112-114. The sheet grows from A1:T110 to
# A1:T114.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 fresh rows right before the old row 108 (whole-row insert
# shifts the existing 108:110 block - and everything below - down by 4
# rows, preserving their formatting/content untouched).
$ws.Rows("108:111").Insert()

# Shared values for every new row.
$mercado  = "Agrícola del Norte S.A. de Arica"
$region   = "Arica y Parinacota"
$fecha    = 45166
$unidad   = '$/bandeja 3 kilos'
$origen   = "Región de Arica y Parinacota"

$rows = @(
    @{ Row = 108; Calidad = "Especial"; Volumen = 250; PrecioMin = 8000; PrecioMax = 9000; PrecioProm = 8400; PrecioKg = 2800 },
    @{ Row = 109; Calidad = "Primera";  Volumen = 300; PrecioMin = 6000; PrecioMax = 7000; PrecioProm = 6500; PrecioKg = 2167 },
    @{ Row = 110; Calidad = "Segunda";  Volumen = 290; PrecioMin = 4000; PrecioMax = 5000; PrecioProm = 4517; PrecioKg = 1506 },
    @{ Row = 111; Calidad = "Tercera";  Volumen = 220; PrecioMin = 2000; PrecioMax = 3000; PrecioProm = 2545; PrecioKg = 848 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value  = 1
    $ws.Cells.Item($row, 2).Value  = $mercado
    $ws.Cells.Item($row, 3).Value  = $region
    $ws.Cells.Item($row, 4).Value  = $fecha
    $ws.Cells.Item($row, 5).Value  = 15
    $ws.Cells.Item($row, 6).Value  = "Fruta"
    $ws.Cells.Item($row, 7).Value  = 100101
    $ws.Cells.Item($row, 8).Value  = "Berries"
    $ws.Cells.Item($row, 9).Value  = 100112025
    $ws.Cells.Item($row, 10).Value = "Frutilla"
    $ws.Cells.Item($row, 11).Value = "Sin especificar"
    $ws.Cells.Item($row, 12).Value = $r.Calidad
    $ws.Cells.Item($row, 13).Value = $r.Volumen
    $ws.Cells.Item($row, 14).Value = $r.PrecioMin
    $ws.Cells.Item($row, 15).Value = $r.PrecioMax
    $ws.Cells.Item($row, 16).Value = $r.PrecioProm
    $ws.Cells.Item($row, 17).Value = $unidad
    $ws.Cells.Item($row, 18).Value = $origen
    $ws.Cells.Item($row, 19).Value = $r.PrecioKg
    $ws.Cells.Item($row, 20).Value = 3
}
